# edit.ps1
# Applies the commit "debug the descriptive data":
#  - "main" sheet: fix R2 count (4676 -> 4671); the describe()-style block in
#    rows 3-5 ("unique"/"top"/"freq") incorrectly included numeric ESG/financial
#    columns (G..W, Y, AC, AG) that should only have the numeric summary stats
#    in rows 6-12 ("mean"/"std"/"min"/"25%"/"50%"/"75%"/"max"). Clear the
#    erroneous categorical-style values from rows 3-5 and fill in the
#    corresponding numeric summary stats in rows 6-12.
#  - "kurtosis" sheet: a missing "beta_5y" variable is inserted at row 16
#    (shifting everything below down by one row), and the trailing
#    "kurtosis_mean" aggregate is recomputed to include it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "main"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("main")

# 1) R2 count: 4676 -> 4671
$ws.Range("R2").Value = 4671

# 2) Columns that wrongly carried "unique"/"top"/"freq" values in rows 3-5.
$statCols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","Y","AC","AG")

foreach ($r in @(3,4,5)) {
    foreach ($col in $statCols) {
        $ws.Range($col + $r).ClearContents()
    }
}

# 3) Numeric summary stats (mean/std/min/25%/50%/75%/max) for rows 6-12.
$rowData = [ordered]@{}
$rowData[6] = @("33.70706578089985","28.33334874496765","27.56358640248537","46.2498465803177","1067.891418394534","1433.696162711349","1361.401587873374","65663.71232334647","29.22150861250581","57.56264489843852","53.91370850003503","0.9796489454155428","0.9913664478159142","6.396770864357124","35.46345637733479","9702984304.326397","33.09201236408795","0.05832860183925946","8.784963576817454","3.503230803844962")
$rowData[7] = @("15.89668564686598","22.83687830675238","18.50104075830211","20.32326413742491","1344.040758952771","1454.214189181862","1222.38798680411","84479.82102692958","72.22063295484895","437.2345687474991","267.7276731218359","0.4824458265927976","0.3360628373851443","8.048556518006695","10.6938924936928","23007356151.29241","11.44163739084912","0.06068279104158179","241.1913262590866","22.9911978011673")
$rowData[8] = @("0.651094295417696","0","0.592213677561109","0.289514866979655","0","0","0.36843861108407","0.2760143558431286","0.5600000000000001","0.320291022","0.6044617289999999","0","0.08530137450098001","-95.689276874","12.2371493409578","195987651.436259","8.965642489855091","8.58034321372881e-05","0.030752545083567","0.01984189847756097")
$rowData[9] = @("21.51932351792345","8.810816534404218","12.2900707911852","30.18440560462735","105.6124499945733","300.942610889746","458.937770674755","9965.260152541205","6.32","7.2156382705","13.25919283227005","0.709206999178785","0.745876159771098","2.0033980695","27.3971368905565","2292567215.601425","24.5297018565522","0.02125709699865266","1.007832366053333","0.5718624931287699")
$rowData[10] = @("31.771844482157","24.4360034294591","23.15046937562665","45.0934205810732","475.7140092095898","968.1687881289047","985.5744766587961","32072.11559044857","12.14","16.83714295776245","24.439861386","0.994373256099218","0.987780430860029","4.7018282818","35.4477774939148","3945383378.846955","32.5367721986685","0.04091676234190243","2.38713078079621","1.108733225396588")
$rowData[11] = @("44.49316748563908","44.69729700948577","39.81119270538638","63.14093354771798","1611.299736103388","2127.211568154549","1900.673544761133","88080.54127036955","27.3","38.014282322","47.0431741845893","1.2655111659106","1.21827959640124","9.3896086345","42.0905114978382","7873099174.596102","40.0118179020734","0.07541939910015569","5.6847155725","2.292691430030193")
$rowData[12] = @("83.5074059365155","93.3945755527551","93.0309717928702","93.3877184731074","7562.459395633622","7296.496723533623","7792.219363268558","582337.7968624146","2050","23572.51592436","11654.5454545451","3.27274160802743","2.15987087037326","87.9587540382","80.4222235607316","405410099022.37","91.5618122495579","1.654364456876971","16123.570794198","1155.515486488235")

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $statCols.Length; $i++) {
        $ws.Range($statCols[$i] + $r).Value = [double]$vals[$i]
    }
}

# ---------------------------------------------------------------------------
# Sheet "kurtosis"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("kurtosis")

# Insert a new row at 16 for "beta_5y", pushing beta_6m.. down by one row.
$ws2.Rows.Item(16).Insert()
$ws2.Range("A16").Value = "beta_5y"
$ws2.Range("B16").Value = [double]"0.914562158242985"

# The trailing "kurtosis_mean" row (now row 40) is recomputed to include the
# newly-added beta_5y kurtosis value.
$ws2.Range("B40").Value = [double]"269.4179111823569"
